$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Land Use Arable"
$ws.Range("A3").Value = "Land Use Fallow"
$ws.Range("A4").Value = "Land Use Perm Past"
$ws.Range("A5").Value = "GHG LUC"
$ws.Range("A6").Value = "GHG Feed"
$ws.Range("A7").Value = "GHG Farm"
$ws.Range("A8").Value = "GHG Processing"
$ws.Range("A9").Value = "GHG Transport"
$ws.Range("A10").Value = "GHG Packging"
$ws.Range("A11").Value = "GHG Retail"
$ws.Range("A12").Value = "Acidification"
$ws.Range("A13").Value = "Eutrophication"
$ws.Range("A14").Value = "Freshwater Withdrawals (FW)"
$ws.Range("A15").Value = "Scarcity-Weighted FW"
